# Scheduled market-data refresh: update currentAveragePrice/LevePrice/LeveProfit
# columns (H-N) on the affected leve rows across the Sheets workbook.
$wb = $excel.ActiveWorkbook

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 6301.1
$ws.Range("I62").Value = 2502
$ws.Range("K62").Value = 2502
$ws.Range("M62").Value = -1878

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 6301.1
$ws.Range("I65").Value = 2502
$ws.Range("K65").Value = 12510
$ws.Range("M65").Value = -9390

# ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3834.25
$ws.Range("I76").Value = 3166.5
$ws.Range("J76").Value = 4502
$ws.Range("K76").Value = 3166.5
$ws.Range("L76").Value = 4502
$ws.Range("M76").Value = -2851.5
$ws.Range("N76").Value = -5132

# ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3834.25
$ws.Range("I79").Value = 3166.5
$ws.Range("J79").Value = 4502
$ws.Range("K79").Value = 3166.5
$ws.Range("L79").Value = 4502
$ws.Range("M79").Value = -2074.5
$ws.Range("N79").Value = -6686

# ALC row 92
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 314
$ws.Range("I92").Value = 314
$ws.Range("K92").Value = 314
$ws.Range("M92").Value = 934

# ALC row 106
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 1679
$ws.Range("I106").Value = 1679
$ws.Range("K106").Value = 1679
$ws.Range("M106").Value = -1048

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1375.8
$ws.Range("I2").Value = 1375.8
$ws.Range("K2").Value = 1375.8
$ws.Range("M2").Value = -1262.8

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3370.4546
$ws.Range("I74").Value = 3143.3333
$ws.Range("J74").Value = 4392.5
$ws.Range("K74").Value = 3143.3333
$ws.Range("L74").Value = 4392.5
$ws.Range("M74").Value = -2269.3333
$ws.Range("N74").Value = -6140.5

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3370.4546
$ws.Range("I77").Value = 3143.3333
$ws.Range("J77").Value = 4392.5
$ws.Range("K77").Value = 15716.6665
$ws.Range("L77").Value = 21962.5
$ws.Range("M77").Value = -11348.6665
$ws.Range("N77").Value = -30698.5

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1375.8
$ws.Range("I116").Value = 1375.8
$ws.Range("K116").Value = 1375.8
$ws.Range("M116").Value = 918.2

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1646.3334
$ws.Range("I122").Value = 1401.6923
$ws.Range("J122").Value = 2282.4
$ws.Range("K122").Value = 4205.0769
$ws.Range("L122").Value = 6847.200000000001
$ws.Range("M122").Value = -1755.0769
$ws.Range("N122").Value = -11747.2

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1930.069
$ws.Range("I132").Value = 1766.4
$ws.Range("K132").Value = 5299.200000000001
$ws.Range("M132").Value = -2769.200000000001

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1375.8
$ws.Range("I3").Value = 1375.8
$ws.Range("K3").Value = 1375.8
$ws.Range("M3").Value = -1261.8

# BSM row 22
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I22").Value = 200
$ws.Range("K22").Value = 200
$ws.Range("M22").Value = -27

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 599.1429000000001
$ws.Range("I94").Value = 599.1429000000001
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 599.1429000000001
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -148.1429000000001
$ws.Range("N94").ClearContents()

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4577.3335
$ws.Range("I107").Value = 1795.5555
$ws.Range("K107").Value = 1795.5555
$ws.Range("M107").Value = 124.4445000000001

# CRP row 64
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 50000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 50000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 50000
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -50496

# CRP row 67
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H67").Value = 50000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 50000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 50000
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -51716

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1990.2593
$ws.Range("I132").Value = 1793.4584
$ws.Range("J132").Value = 3564.6667
$ws.Range("K132").Value = 5380.3752
$ws.Range("L132").Value = 10694.0001
$ws.Range("M132").Value = -2850.3752
$ws.Range("N132").Value = -15754.0001

# CUL row 6
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 444.92856
$ws.Range("I6").Value = 125.44444
$ws.Range("K6").Value = 376.33332
$ws.Range("M6").Value = -263.33332

# CUL row 12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 232.5
$ws.Range("J12").Value = 259.2
$ws.Range("L12").Value = 777.5999999999999
$ws.Range("N12").Value = -1123.6

# CUL row 34
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1909.8
$ws.Range("J34").Value = 2254.52
$ws.Range("L34").Value = 6763.559999999999
$ws.Range("N34").Value = -6931.559999999999

# CUL row 92
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 3612.125
$ws.Range("J92").Value = 5019.6
$ws.Range("L92").Value = 15058.8
$ws.Range("N92").Value = -17554.8

# CUL row 129
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2051.111
$ws.Range("J129").Value = 2932
$ws.Range("L129").Value = 8796
$ws.Range("N129").Value = -18796

# CUL row 133
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 8677.786
$ws.Range("I113").Value = 7686.125
$ws.Range("K113").Value = 7686.125
$ws.Range("M113").Value = -5516.125

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1479.8
$ws.Range("I122").Value = 863.6667
$ws.Range("K122").Value = 2591.0001
$ws.Range("M122").Value = -141.0001000000002

# LTW row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2969.2666
$ws.Range("I82").Value = 1171.1111
$ws.Range("J82").Value = 5666.5
$ws.Range("K82").Value = 1171.1111
$ws.Range("L82").Value = 5666.5
$ws.Range("M82").Value = -810.1111000000001
$ws.Range("N82").Value = -6388.5

# LTW row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2969.2666
$ws.Range("I85").Value = 1171.1111
$ws.Range("J85").Value = 5666.5
$ws.Range("K85").Value = 1171.1111
$ws.Range("L85").Value = 5666.5
$ws.Range("M85").Value = 76.88889999999992
$ws.Range("N85").Value = -8162.5

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4700.8
$ws.Range("I132").Value = 5676
$ws.Range("K132").Value = 17028
$ws.Range("M132").Value = -14498

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3175.6924
$ws.Range("I122").Value = 2062.375
$ws.Range("J122").Value = 4957
$ws.Range("K122").Value = 6187.125
$ws.Range("L122").Value = 14871
$ws.Range("M122").Value = -3737.125
$ws.Range("N122").Value = -19771

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3861.6
$ws.Range("I126").Value = 1540.1666
$ws.Range("K126").Value = 4620.4998
$ws.Range("M126").Value = -2150.4998
